$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '57.872.03'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -5.56%  '
$c.Style = "Normal"

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.899.63'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -3.81%  '
$c.Style = "Normal"

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '550.52'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -3.78%  '
$c.Style = "Normal"

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '122.59'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -4.87%  '
$c.Style = "Normal"

# Row 7
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '2.894.09'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -3.88%  '
$c.Style = "Normal"

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.494'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -0.79%  '
$c.Style = "Normal"

# Row 10
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -7.72%  '
$c.Style = "Normal"

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '4.72'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -9.18%  '
$c.Style = "Normal"

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.434'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +0.42%  '
$c.Style = "Normal"

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000212'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -6.51%  '
$c.Style = "Normal"

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '31.82'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -4.48%  '
$c.Style = "Normal"

# Row 15
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.369.78'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -4.05%  '
$c.Style = "Normal"

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.886.11'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -4.17%  '
$c.Style = "Normal"

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '57.737.17'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -5.89%  '
$c.Style = "Normal"

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.49'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +3.36%  '
$c.Style = "Normal"

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '408.35'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -6.81%  '
$c.Style = "Normal"

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.93'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -2.32%  '
$c.Style = "Normal"

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.657'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -1.20%  '
$c.Style = "Normal"

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.77'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -5.99%  '
$c.Style = "Normal"

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '12.62'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.Style = "Normal"

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '77.15'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c.Style = "Normal"

# Row 26
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.Style = "Normal"

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -2.29%  '
$c.Style = "Normal"

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.24'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -1.37%  '
$c.Style = "Normal"

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.91'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -3.17%  '
$c.Style = "Normal"

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -3.09%  '
$c.Style = "Normal"

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '24.67'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -3.67%  '
$c.Style = "Normal"

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0957'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +1.41%  '
$c.Style = "Normal"

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.912'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -4.93%  '
$c.Style = "Normal"

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.03'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -11.38%  '
$c.Style = "Normal"

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.36'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -4.49%  '
$c.Style = "Normal"

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '48.31'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -3.86%  '
$c.Style = "Normal"

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '8.49'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +8.89%  '
$c.Style = "Normal"

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0₃0623'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -8.95%  '
$c.Style = "Normal"

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0345'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -6.03%  '
$c.Style = "Normal"

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.106'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -3.43%  '
$c.Style = "Normal"

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.611.84'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -1.75%  '
$c.Style = "Normal"

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '360.18'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -3.78%  '
$c.Style = "Normal"

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -3.81%  '
$c.Style = "Normal"

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.229'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -3.08%  '
$c.Style = "Normal"

# Row 47
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '117.82'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -3.05%  '
$c.Style = "Normal"

# Row 48
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -0.82%  '
$c.Style = "Normal"

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.94'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -1.97%  '
$c.Style = "Normal"

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '22.72'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -3.53%  '
$c.Style = "Normal"

# Row 51
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.95'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -4.04%  '
$c.Style = "Normal"
